$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Simplistic emissions model"
$ws.Range("C5").Value = "For each working condition, emissions for CO2,CO,Nox,THC is estimated"

$ws.Range("A6").Select()
